# Generate Report for Handback
# Update the "Correspond Handback Datetime" / completion timestamps for the
# most recent handback entry (b1e3693c-... file, row 2) on both the zh-cn
# and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 08:53:23"
$wsZhCn.Range("H2").Value = "2016-03-23 08:53:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 08:53:27"
$wsDeDe.Range("H2").Value = "2016-03-23 08:53:54"
